$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) from 45243 to 45244 for rows 2 through 9
for ($row = 2; $row -le 9; $row++) {
    $cell = $ws.Range("C$row")
    if ($cell.Value2 -eq 45243) {
        $cell.Value = 45244
    }
}
